$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$headers = @("g", "n", "sqgrupos", "sqerros", "dfgrupos", "dferros", "sqtotal", "dftotal", "msqgrupos", "msqerros", "f", "sig")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2  # headers start at column B (2)
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Apply the existing header style (already used on B1:G1) to the newly
# added header cells H1:M1 so every header cell is formatted consistently.
$ws.Range("B1:G1").Copy()
$ws.Range("H1:M1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (rows 2-6), columns A..M ---
$data = @(
    @(0, 4, 4, 4,    4,   3, 12, 8,    15, 1.333,   0.333, 4.003,  0.035),
    @(1, 4, 6, 1000, 100, 3, 20, 1100, 23, 333.333, 5,     66.667, 0),
    @(2, 4, 5, 500,  480, 3, 16, 980,  19, 166.667, 30,    5.556,  0.008),
    @(3, 5, 5, 180,  170, 4, 20, 350,  24, 45,      8.5,   5.294,  0.004),
    @(4, 4, 4, 200,  120, 3, 12, 320,  15, 66.667,  10,    6.667,  0.007)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$c]
    }
}

# Apply the existing "A column" style (used on A2:A5) to the newly added
# A6 cell so the whole A column stays formatted consistently.
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
